$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2031746031746032
$ws.Range("C2").Value = 0.546031746031746
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.1079365079365079
$ws.Range("B3").Value = 0.005714285714285714
$ws.Range("C3").Value = 0.02857142857142857
$ws.Range("J3").Value = 0.005714285714285714
$ws.Range("P3").Value = 0.7657142857142857
$ws.Range("S3").Value = 0.1942857142857143
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.04958677685950413
$ws.Range("D6").Value = 0.01239669421487603
$ws.Range("E6").Value = 0.004132231404958678
$ws.Range("F6").Value = 0.07024793388429752
$ws.Range("J6").Value = 0.2520661157024793
$ws.Range("O6").Value = 0.02066115702479339
$ws.Range("Q6").Value = 0.1818181818181818
$ws.Range("R6").Value = 0.0743801652892562
$ws.Range("S6").Value = 0.3347107438016529
$ws.Range("B7").Value = 0.09743589743589744
$ws.Range("D7").Value = 0.02564102564102564
$ws.Range("E7").Value = 0.005128205128205128
$ws.Range("F7").Value = 0.08717948717948718
$ws.Range("J7").Value = 0.08717948717948718
$ws.Range("O7").Value = 0.01538461538461539
$ws.Range("Q7").Value = 0.2564102564102564
$ws.Range("R7").Value = 0.04102564102564103
$ws.Range("S7").Value = 0.3846153846153846
$ws.Range("B8").Value = 0.1165755919854281
$ws.Range("D8").Value = 0.02003642987249545
$ws.Range("F8").Value = 0.0546448087431694
$ws.Range("J8").Value = 0.09836065573770492
$ws.Range("O8").Value = 0.01821493624772313
$ws.Range("Q8").Value = 0.1493624772313297
$ws.Range("R8").Value = 0.06557377049180328
$ws.Range("S8").Value = 0.4772313296903461
$ws.Range("B9").Value = 0.07518796992481203
$ws.Range("D9").Value = 0.01503759398496241
$ws.Range("F9").Value = 0.06015037593984962
$ws.Range("J9").Value = 0.05263157894736842
$ws.Range("O9").Value = 0.02631578947368421
$ws.Range("Q9").Value = 0.1804511278195489
$ws.Range("R9").Value = 0.06390977443609022
$ws.Range("S9").Value = 0.5263157894736842
$ws.Range("B10").Value = 0.115695067264574
$ws.Range("D10").Value = 0.01704035874439462
$ws.Range("F10").Value = 0.08609865470852018
$ws.Range("J10").Value = 0.08699551569506726
$ws.Range("O10").Value = 0.01524663677130045
$ws.Range("Q10").Value = 0.1901345291479821
$ws.Range("R10").Value = 0.08071748878923767
$ws.Range("S10").Value = 0.4080717488789238
$ws.Range("G11").Value = 0.1409395973154362
$ws.Range("J11").Value = 0.05704697986577181
$ws.Range("K11").Value = 0.1912751677852349
$ws.Range("L11").Value = 0.587248322147651
$ws.Range("S11").Value = 0.02348993288590604
$ws.Range("G12").Value = 0.7431693989071039
$ws.Range("J12").Value = 0.1912568306010929
$ws.Range("K12").Value = 0.01092896174863388
$ws.Range("L12").Value = 0.02185792349726776
$ws.Range("S12").Value = 0.03278688524590164
$ws.Range("F15").Value = 0.045
$ws.Range("H15").Value = 0.15
$ws.Range("I15").Value = 0.105
$ws.Range("J15").Value = 0.26
$ws.Range("K15").Value = 0.06
$ws.Range("M15").Value = 0.015
$ws.Range("O15").Value = 0.055
$ws.Range("S15").Value = 0.31
$ws.Range("F16").Value = 0.015
$ws.Range("H16").Value = 0.215
$ws.Range("I16").Value = 0.105
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.075
$ws.Range("M16").Value = 0.01
$ws.Range("O16").Value = 0.025
$ws.Range("S16").Value = 0.155
$ws.Range("F17").Value = 0.01601830663615561
$ws.Range("H17").Value = 0.2196796338672769
$ws.Range("I17").Value = 0.1167048054919908
$ws.Range("J17").Value = 0.4050343249427917
$ws.Range("K17").Value = 0.06636155606407322
$ws.Range("M17").Value = 0.01601830663615561
$ws.Range("N17").Value = 0.002288329519450801
$ws.Range("O17").Value = 0.05949656750572083
$ws.Range("S17").Value = 0.09839816933638444
$ws.Range("F18").Value = 0.02395209580838323
$ws.Range("H18").Value = 0.1856287425149701
$ws.Range("I18").Value = 0.1077844311377246
$ws.Range("J18").Value = 0.3712574850299401
$ws.Range("K18").Value = 0.1077844311377246
$ws.Range("M18").Value = 0.02395209580838323
$ws.Range("N18").Value = 0.005988023952095809
$ws.Range("O18").Value = 0.04191616766467066
$ws.Range("S18").Value = 0.1317365269461078
$ws.Range("F19").Value = 0.01193820224719101
$ws.Range("H19").Value = 0.2429775280898877
$ws.Range("I19").Value = 0.1102528089887641
$ws.Range("J19").Value = 0.3223314606741573
$ws.Range("K19").Value = 0.1165730337078652
$ws.Range("M19").Value = 0.01264044943820225
$ws.Range("O19").Value = 0.05617977528089887
$ws.Range("S19").Value = 0.1271067415730337
